# Regenerate orders with updated distance/size codes.
# The experiment's Distance codes (D64->D69, D51->D55, D80->D86) and the
# Size code S30->S31 changed; every string built from them (Condition,
# Filename_Left, Filename_Right, Distance, Size lookup lists) needs the
# same substitution applied wherever it occurs as a substring.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

$null = $used.Replace("D64", "D69")
$null = $used.Replace("D51", "D55")
$null = $used.Replace("D80", "D86")
$null = $used.Replace("S30", "S31")
